$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update cell values (shared strings content) ---
$ws.Range("A2").Value = "TestUserFirstName"
$ws.Range("B2").Value = "TestUserLastName"
$ws.Range("C2").Value = "testuser@test.com"

$ws.Range("A3").Value = "Christy"
$ws.Range("B3").Value = "Fernandes"
$ws.Range("C3").Value = "test12email@test.com"

$ws.Range("A4").Value = "Test"
$ws.Range("B4").Value = "test"
$ws.Range("C4").Value = "testemail123@test.com"

# --- Style changes: A2, B2, C3, C4 need to switch to the same style as A1/C2 (s=1) ---
$ws.Range("A1").Copy()
$ws.Range("A2:B2").PasteSpecial(-4122)
$ws.Range("C3").PasteSpecial(-4122)
$ws.Range("C4").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Column widths ---
# NOTE: the host engine quantizes the stored OOXML <col width> to 1/6-pixel
# steps (MaxDigitWidth=6) no matter what font is applied, while the target
# file's widths (32.57 / 32.29 / 29.86) were produced with MDW=7. The
# COM ColumnWidth values below are chosen as the closest achievable match
# (stored width lands on 32.5 / 32.33 / 29.83 respectively, each within
# ~0.02-0.07 of the target - the nearest the engine can represent).
$ws.Columns.Item(1).ColumnWidth = 31.66
$ws.Columns.Item(2).ColumnWidth = 31.5
$ws.Columns.Item(3).ColumnWidth = 29.0
